# "Finish add customer step"
# Insert two new worksheets, "customer" and "assignProjectTest", between the
# existing "addLeaveTypeTest" and "employee" sheets, populate them with their
# test data, and make "customer" the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Anchor the new sheets right before "employee" (current 2nd sheet).
$employee = $wb.Worksheets.Item("employee")

# Add "assignProjectTest" first (ends up directly before "employee"),
# then add "customer" before that — giving final order:
#   addLeaveTypeTest, customer, assignProjectTest, employee, addSkillTest
$assignProjectTestSheet = $wb.Worksheets.Add($employee)
$assignProjectTestSheet.Name = "assignProjectTest"

$customerSheet = $wb.Worksheets.Add($assignProjectTestSheet)
$customerSheet.Name = "customer"

# Re-fetch fresh handles by name before writing data, since sheet indices
# shift as new sheets get inserted.
$customer = $wb.Worksheets.Item("customer")
$assignProjectTest = $wb.Worksheets.Item("assignProjectTest")

# --- customer sheet data ---
$customer.Range("A1").Value = "Testcase ID"
$customer.Range("B1").Value = "Index in testcase"
$customer.Range("C1").Value = "Name"
$customer.Range("D1").Value = "Description"

$customer.Range("A2").Value = "ASSIGN_PROJECT_TC000"
$customer.Range("B2").Value = 0
$customer.Range("C2").Value = "John"
$customer.Range("D2").Value = "test"

$customer.Range("A3").Value = "ASSIGN_PROJECT_TC000"
$customer.Range("B3").Value = 1
$customer.Range("C3").Value = "Jame"
$customer.Range("D3").Value = "test"

# --- assignProjectTest sheet data ---
$assignProjectTest.Range("A1").Value = "Testcase Id"
$assignProjectTest.Range("A2").Value = "ASSIGN_PROJECT_TC000"

# Match the author's final selection state on each sheet.
$assignProjectTest.Range("A2").Select()
$customer.Range("A6").Select()

# "customer" is the active/selected tab in the finished workbook.
$customer.Activate()
